$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "... en y ajoutant un identifiant en debut (<< server >>) et fin
#    (<< end >>) de chaine ..." becomes "... identifiant (<< -1 >>) en
#    debut et fin de chaine ..."
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    " en début (« server ») et fin (« end ») de chaîne afin de permettre aux robots de différencier les messages provenant du serveur des messages provenant des autres robots.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " (« -1 ») en début et fin de chaîne afin de permettre aux robots de différencier les messages provenant du serveur des messages provenant des autres robots.",
    1) | Out-Null

# ---------------------------------------------------------------------
# 2) "... ainsi que des ID, position et vitesses ..." becomes
#    "... ainsi que des IDs, positions et vitesses ..."
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "ainsi que des ID, position et vitesses",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "ainsi que des IDs, positions et vitesses",
    1) | Out-Null

# ---------------------------------------------------------------------
# 3) The quoted example "<< server\n0\n0\n255\n[...]\nend >>" becomes
#    "<< -1\n0\n0\n255\n[...]\n-1 >>"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "« server\n0\n0\n255\n[…]\nend »",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "« -1\n0\n0\n255\n[…]\n-1 »",
    1) | Out-Null

# ---------------------------------------------------------------------
# 4) Move the "_GoBack" bookmark: it used to sit after "On prend ici le
#    temps serveur," -- now it should sit right after "positions" (the
#    point of the last text edit), i.e. just before " et vitesses".
#    Adding a bookmark named "_GoBack" again relocates the existing one.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(
    "ainsi que des IDs, positions",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rng) | Out-Null

# ---------------------------------------------------------------------
# 5) Second table: "server" -> "-1", "end" -> "-1", center every cell's
#    text, and shrink the first/last grid columns.
# ---------------------------------------------------------------------
$t2 = $d.Tables.Item(2)

$t2.Cell(1, 1).Range.Text = "-1"
$t2.Cell(1, 11).Range.Text = "-1"

for ($r = 1; $r -le $t2.Rows.Count; $r++) {
    for ($c = 1; $c -le $t2.Columns.Count; $c++) {
        $t2.Cell($r, $c).Range.ParagraphFormat.Alignment = 1
    }
}

$t2.Columns.Item(1).Width = 766 / 20.0
$t2.Columns.Item(11).Width = 484 / 20.0
